$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Matches_A" column (B) and "Matches_B" column (H) values
$ws.Range("B2").Value = 5
$ws.Range("H2").Value = 12

$ws.Range("B3").Value = 8
$ws.Range("H3").Value = 8

$ws.Range("B4").Value = 9
$ws.Range("H4").Value = 8

$ws.Range("B5").Value = 10
$ws.Range("H5").Value = 14
